$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new Airbyte metadata columns ---------------------------
# Original layout:  A=pca_code  B=pca_iden  C=pca_labe  D=updated_at
# Target layout:    A=_airbyte_ab_id  B=_airbyte_emitted_at  C=pca_code
#                   D=pca_iden  E=pca_labe  F=_airbyte_additional_properties
#                   G=source_file_path  H=updated_at

# Insert 2 columns before the current A (becomes A,B new / C,D,E,F old data)
$ws.Range("A:B").Insert()
# Insert 2 more columns before what is now column F (old "updated_at" column)
$ws.Range("F:G").Insert()

# New A1/B1 did not inherit the bold/bordered header style (they were the
# leftmost columns before insertion), so copy it over from a neighbouring
# header cell that already carries it. F1/G1 picked it up automatically.
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("F1").Value = "_airbyte_additional_properties"
$ws.Range("G1").Value = "source_file_path"

# --- Per-row data ----------------------------------------------------------
$ids = @(
    "939401d6-7437-40a4-9842-c496c3b987aa",
    "9702a27d-6ff5-4164-b9f4-7e9e67e3189b",
    "4f41b759-02be-47a5-a80a-9c4d9e311969",
    "fca49d48-b0fd-4208-871b-f5b65f316945",
    "7358340c-ebbb-4124-81bc-ab70b161fe92",
    "c5e510aa-9ed6-42ed-826b-3f848e23c111",
    "d0a1c6ee-8f05-4dd0-a35d-09bef97c1889",
    "0bba6200-856e-4b2a-9c02-5603edbead62",
    "aebd5cbd-19b9-4456-9be6-b4a9b965c31a",
    "1b37a342-dd50-4d2f-9582-e38165f01efe",
    "86af1fa1-84d3-4cbd-9cfa-439c0ba6dd00"
)

$sourceFile = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/PAYMENT_CHANNEL/2024_08_06_1722929004063_0.parquet"
$emittedAt = 45510.3079196875
$updatedAt = 45511.29525109382

for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value = $ids[$i]

    $ws.Cells.Item($row, 2).Value = $emittedAt
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 7).Value = $sourceFile

    $ws.Cells.Item($row, 8).Value = $updatedAt
    $ws.Cells.Item($row, 8).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
